$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.301.22"
$ws.Range("E2").Value = "  +2.34%  "
$ws.Range("D3").Value = "2.539.23"
$ws.Range("E3").Value = "  +3.08%  "
$ws.Range("D5").Formula = "'582.93"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").Formula = "'153.09"
$ws.Range("E6").Value = "  +4.96%  "
$ws.Range("D7").Formula = "'1.00"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Formula = "'0.539"
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").Formula = "'0.356"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Formula = "'29.75"
$ws.Range("E13").Value = "  +2.59%  "
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").Value = "2.993.70"
$ws.Range("E15").Value = "  +2.87%  "
$ws.Range("D16").Value = "63.941.07"
$ws.Range("E16").Value = "  +1.74%  "
$ws.Range("D17").Value = "2.532.92"
$ws.Range("E17").Value = "  +2.60%  "
$ws.Range("D18").Formula = "'7.93"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("D19").Formula = "'11.02"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("E20").Value = "  +3.49%  "
$ws.Range("D21").Formula = "'328.05"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Formula = "'2.26"
$ws.Range("E22").Value = "  +1.44%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Formula = "'10.11"
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("D25").Formula = "'65.63"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("D26").Formula = "'662.18"
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("D27").Formula = "'0.0000104"
$ws.Range("E27").Value = "  +5.59%  "
$ws.Range("D28").Value = "2.685.13"
$ws.Range("E28").Value = "  +3.82%  "
$ws.Range("E29").Value = "  +2.72%  "
$ws.Range("D30").Formula = "'0.998"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Formula = "'8.08"
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").Formula = "'1.87"
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("E33").Value = "  +1.81%  "
$ws.Range("D34").Formula = "'0.999"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  +1.45%  "
$ws.Range("E36").Value = "  +1.78%  "
$ws.Range("E37").Value = "  +3.55%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Formula = "'2.83"
$ws.Range("E38").Value = "  +2.58%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Formula = "'0.372"
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").Formula = "'18.93"
$ws.Range("E40").Value = "  +1.14%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Formula = "'152.01"
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("E42").Value = "  +3.33%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Formula = "'159.33"
$ws.Range("E44").Value = "  +3.21%  "
$ws.Range("E45").Value = "  -3.11%  "
$ws.Range("E46").Value = "  +1.44%  "
$ws.Range("E47").Value = "  +1.99%  "
$ws.Range("D48").Formula = "'21.10"
$ws.Range("E48").Value = "  +3.92%  "
$ws.Range("E49").Value = "  +2.26%  "
$ws.Range("D50").Formula = "'0.0523"
$ws.Range("E50").Value = "  +2.09%  "
$ws.Range("D51").Formula = "'0.0230"
$ws.Range("E51").Value = "  +2.31%  "
